$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text: volume number and report week dates ---
$ws.Range("A8").Value = "Volume 30   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/10/2023  Through  7/16/2023"

# --- Update crime statistics table (rows 14-29) ---

# Row 14
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "0"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "***.*"
$ws.Range("N14").Value = -92.857142857142

# Row 15
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -33.333333333333
$ws.Range("J15").Value = 14
$ws.Range("K15").Value = 0
$ws.Range("M15").Value = -6.666666666666

# Row 16
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 11
$ws.Range("H16").Value = -8.333333333333
$ws.Range("I16").Value = 99
$ws.Range("J16").Value = 83
$ws.Range("K16").Value = 19.277108433734
$ws.Range("L16").Value = 15.116279069767
$ws.Range("M16").Value = -39.263803680981
$ws.Range("N16").Value = -84.792626728110

# Row 17
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -12.5
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 33
$ws.Range("H17").Value = -27.272727272727
$ws.Range("I17").Value = 174
$ws.Range("J17").Value = 191
$ws.Range("K17").Value = -8.900523560209
$ws.Range("L17").Value = 6.748466257668
$ws.Range("M17").Value = 5.454545454545
$ws.Range("N17").Value = -62.419006479481

# Row 18
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -60
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 91
$ws.Range("J18").Value = 111
$ws.Range("K18").Value = -18.018018018018
$ws.Range("L18").Value = 7.058823529411
$ws.Range("M18").Value = -14.953271028037
$ws.Range("N18").Value = -77.25

# Row 19
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -20
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 29
$ws.Range("H19").Value = 10.344827586206
$ws.Range("I19").Value = 187
$ws.Range("J19").Value = 211
$ws.Range("K19").Value = -11.374407582938
$ws.Range("L19").Value = 1.081081081081
$ws.Range("M19").Value = 29.861111111111
$ws.Range("N19").Value = -0.531914893617

# Row 20
$ws.Range("C20").Value = 5
$ws.Range("E20").Value = 66.666666666666
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 25
$ws.Range("I20").Value = 54
$ws.Range("J20").Value = 76
$ws.Range("K20").Value = -28.947368421052
$ws.Range("L20").Value = 12.5
$ws.Range("M20").Value = 12.5
$ws.Range("N20").Value = -84.438040345821

# Row 21
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = -20
$ws.Range("F21").Value = 90
$ws.Range("H21").Value = -6.25
$ws.Range("I21").Value = 620
$ws.Range("J21").Value = 691
$ws.Range("K21").Value = -10.274963820549
$ws.Range("L21").Value = 6.529209621993
$ws.Range("M21").Value = -5.053598774885
$ws.Range("N21").Value = -70.560303893637

# Row 22
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "0"
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "***.*"
$ws.Range("L22").Value = 11.111111111111
$ws.Range("M22").Value = -33.333333333333

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "***.*"
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 10
$ws.Range("H23").Value = -20
$ws.Range("I23").Value = 48
$ws.Range("K23").Value = 6.666666666666
$ws.Range("L23").Value = -11.111111111111
$ws.Range("M23").Value = 4.347826086956

# Row 24
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = -13.333333333333
$ws.Range("I24").Value = 437
$ws.Range("J24").Value = 388
$ws.Range("K24").Value = 12.628865979381
$ws.Range("L24").Value = 45.182724252491
$ws.Range("M24").Value = 14.099216710182

# Row 25
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 61
$ws.Range("G25").Value = 32
$ws.Range("H25").Value = 90.625
$ws.Range("I25").Value = 292
$ws.Range("J25").Value = 220
$ws.Range("K25").Value = 32.727272727272
$ws.Range("L25").Value = 64.971751412429
$ws.Range("M25").Value = -31.132075471698

# Row 26
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"
$ws.Range("E26").Value = -100
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 19
$ws.Range("K26").Value = 5.263157894736
$ws.Range("L26").Value = 53.846153846153

# Row 27
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 13
$ws.Range("J27").Value = 18
$ws.Range("K27").Value = -27.777777777777
$ws.Range("L27").Value = -55.172413793103

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("F28").NumberFormat = "@"
$ws.Range("F28").Value = "0"
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = -100
$ws.Range("M28").Value = -77.419354838709
$ws.Range("N28").Value = -93.137254901960

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$ws.Range("F29").NumberFormat = "@"
$ws.Range("F29").Value = "0"
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = -100
$ws.Range("M29").Value = -76.923076923076
$ws.Range("N29").Value = -93.478260869565
